$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tube")

# --- Row 1 header changes (O/P/Q column headers reassigned) ---
$ws.Range("O1").Value = "Splash"
$ws.Range("P1").Value = "cooltime"
$ws.Range("Q1").Value = "abnormal_value"

# --- Row 5: K5 content change + P5 -> O5 move ---
$ws.Range("K5").Value = "{(attack01, 0, 0)}"
$ws.Range("O5").Value = 0
$ws.Range("P5").ClearContents()

# --- Row 6: P6 -> O6 move ---
$ws.Range("O6").Value = 1
$ws.Range("P6").ClearContents()

# --- Row 7: P7 -> O7 move ---
$ws.Range("O7").Value = 0
$ws.Range("P7").ClearContents()

# --- Row 8: O8 (slow_01) removed, P8 -> O8 move ---
$ws.Range("O8").Value = 0
$ws.Range("P8").ClearContents()

# --- Row 9: O9 (snare_04) removed, P9 -> O9 move ---
$ws.Range("O9").Value = 0
$ws.Range("P9").ClearContents()

# --- Row 10: Q10 -> P10 move ---
$ws.Range("P10").Value = 2
$ws.Range("Q10").ClearContents()

# --- Row 11: Q11 -> P11 move ---
$ws.Range("P11").Value = 2
$ws.Range("Q11").ClearContents()

# --- Row 12: Q12 -> P12 move ---
$ws.Range("P12").Value = 3
$ws.Range("Q12").ClearContents()

# --- New rows 13-23 ---
# Row 13
$ws.Range("A13").Value = 4100
$ws.Range("B13").Value = "Relic_01"
$ws.Range("C13").Value = "슬로우_01"
$ws.Range("D13").Value = "relic"
$ws.Range("E13").Value = "Weakness"
$ws.Range("F13").Value = "C"
$ws.Range("Q13").Value = "slow_01"

# Row 14
$ws.Range("A14").Value = 4101
$ws.Range("B14").Value = "Relic_02"
$ws.Range("C14").Value = "스네어_01"
$ws.Range("D14").Value = "relic"
$ws.Range("E14").Value = "Weakness"
$ws.Range("F14").Value = "C"
$ws.Range("Q14").Value = "snare_01"

# Row 15
$ws.Range("A15").Value = 5100
$ws.Range("B15").Value = "no107_style"
$ws.Range("C15").Value = "107식 스타일"
$ws.Range("D15").Value = "style"
$ws.Range("E15").Value = "Weakness"
$ws.Range("F15").Value = "A"
$ws.Range("G15").Value = "no107_skill"
$ws.Range("H15").Value = "{(0.7)}"
$ws.Range("I15").Value = "melee"
$ws.Range("J15").Value = 10

# Row 16
$ws.Range("A16").Value = 5101
$ws.Range("B16").Value = "no107_enhancer"
$ws.Range("C16").Value = "107식 인핸서"
$ws.Range("D16").Value = "enhancer"
$ws.Range("E16").Value = "Weakness"
$ws.Range("F16").Value = "A"
$ws.Range("K16").Value = "{(no107_skill, 0, 1)}"
$ws.Range("O16").Value = 0

# Row 17
$ws.Range("A17").Value = 5102
$ws.Range("B17").Value = "no107_cooler"
$ws.Range("C17").Value = "107식 쿨러"
$ws.Range("D17").Value = "cooler"
$ws.Range("E17").Value = "Weakness"
$ws.Range("F17").Value = "A"
$ws.Range("P17").Value = 1.5

# Row 18
$ws.Range("A18").Value = 5103
$ws.Range("B18").Value = "no108_style"
$ws.Range("C18").Value = "108식 스타일"
$ws.Range("D18").Value = "style"
$ws.Range("E18").Value = "Weakness"
$ws.Range("F18").Value = "A"
$ws.Range("G18").Value = "no107_attack1"
$ws.Range("H18").Value = "{(0.7)}"
$ws.Range("I18").Value = "melee"
$ws.Range("J18").Value = 10

# Row 19
$ws.Range("A19").Value = 5104
$ws.Range("B19").Value = "no108_enhancer"
$ws.Range("C19").Value = "108식 인행서"
$ws.Range("D19").Value = "enhancer"
$ws.Range("E19").Value = "Weakness"
$ws.Range("F19").Value = "A"
$ws.Range("K19").Value = "{(no107_attack1, 0, 2)}"
$ws.Range("O19").Value = 0

# Row 20
$ws.Range("A20").Value = 5105
$ws.Range("B20").Value = "no108_cooler"
$ws.Range("C20").Value = "108식 쿨러"
$ws.Range("D20").Value = "cooler"
$ws.Range("E20").Value = "Weakness"
$ws.Range("F20").Value = "A"
$ws.Range("P20").Value = 1.5

# Row 21
$ws.Range("A21").Value = 5106
$ws.Range("B21").Value = "jake_style"
$ws.Range("C21").Value = "제이크 스타일"
$ws.Range("D21").Value = "style"
$ws.Range("E21").Value = "Weakness"
$ws.Range("F21").Value = "A"
$ws.Range("G21").Value = "jake_attack"
$ws.Range("H21").Value = "none"
$ws.Range("I21").Value = "bounce"
$ws.Range("J21").Value = 10

# Row 22
$ws.Range("A22").Value = 5107
$ws.Range("B22").Value = "jake_enhancer"
$ws.Range("C22").Value = "제이크 인핸서"
$ws.Range("D22").Value = "enhancer"
$ws.Range("E22").Value = "Weakness"
$ws.Range("F22").Value = "A"
$ws.Range("M22").Value = "{(jake_attack, 4, 4)}"
$ws.Range("O22").Value = 1

# Row 23
$ws.Range("A23").Value = 5108
$ws.Range("B23").Value = "jake_cooler"
$ws.Range("C23").Value = "제이크 쿨러"
$ws.Range("D23").Value = "cooler"
$ws.Range("E23").Value = "Weakness"
$ws.Range("F23").Value = "A"
$ws.Range("P23").Value = 3

# --- Update selection to match target (B21) ---
$ws.Range("B21").Select()